$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Structural edits: insert two new rows (new rows 4 & 5, pushing the old
#    "Units" table down) and a new column C (pushing the old "Description"
#    column from C to D).
# ---------------------------------------------------------------------------
$ws.Rows("4:5").Insert()
$ws.Columns("C").Insert()

# ---------------------------------------------------------------------------
# 2. New row 4: "Default From Row" = 8  (A4/B4), formatted like the other
#    label rows (style of A1/A2/A3).
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("C1").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("C1").Copy()
$ws.Range("D4").PasteSpecial(-4122)

$ws.Range("A4").Value = "Default From Row"
$ws.Range("B4").Value = 8

# ---------------------------------------------------------------------------
# 3. New row 5: "Default To Row" = 10, with explanatory text in D5.
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("C1").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("D1").Copy()
$ws.Range("D5").PasteSpecial(-4122)

$ws.Range("A5").Value = "Default To Row"
$ws.Range("B5").Value = 10
$ws.Range("D5").Value = 'Use this to select the row with the default "to" conversion value when the page loads (default is 8)'

# ---------------------------------------------------------------------------
# 4. New column C header ("Description") on row 7, matching the bold style
#    already used by A7/B7 ("Name" / "Factor").
# ---------------------------------------------------------------------------
$ws.Range("B7").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("C7").Value = "Description"

# ---------------------------------------------------------------------------
# 5. New column C values (rows 8-17): human readable unit names. Build the
#    font/style used by these new cells once, then stamp it onto the range
#    via PasteSpecial so it is only created a single time.
# ---------------------------------------------------------------------------
$ws.Range("ZZ1").Value = "x"
$ws.Range("ZZ1").Font.Name = "Calibri"
$ws.Range("ZZ1").Font.ThemeFont = 1
$ws.Range("ZZ1").Font.Size = 11
$ws.Range("ZZ1").Font.ThemeColor = 1
$ws.Range("ZZ1").Copy()
$ws.Range("C8:C17").PasteSpecial(-4122)
$ws.Range("ZZ1").Clear()

$ws.Range("C8").Value = "square meter"
$ws.Range("C9").Value = "square millimeter"
$ws.Range("C10").Value = "square foot"
$ws.Range("C11").Value = "square inch"
$ws.Range("C12").Value = "square centimeter"
$ws.Range("C13").Value = "square mile"
$ws.Range("C14").Value = "square survey mile"
$ws.Range("C15").Value = "square kilometer"
$ws.Range("C16").Value = "square yard"
$ws.Range("C17").Value = "acre"

# ---------------------------------------------------------------------------
# 6. Row height tweaks to match the refreshed layout.
# ---------------------------------------------------------------------------
$ws.Rows(6).RowHeight = 16.5
$ws.Rows(9).RowHeight = 51
$ws.Rows(10).RowHeight = 25.5
$ws.Range("11:17").RowHeight = 15

$excel.CutCopyMode = 0
